$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the MovieID sequence: column E (5th movie) should be 4, not 2
$ws.Range("E1").Value = 4

# Update the active selection to match the authored state
$ws.Range("G13").Select()
